# Apply the target edit to the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the values for the joint-probability table (rows 5-8)
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 3

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0.1
$ws.Range("D6").Value = 0.1
$ws.Range("E6").Value = 0.3

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 0.05
$ws.Range("D7").Value = 0.15
$ws.Range("E7").Value = 0.2

$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 0.1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

# Add the new notes / text in column K, plus the new E(Z|y=2) example with its formula
$ws.Range("K14").Value = "This time we even have Conditional Probability:"
$ws.Range("K16").Value = "And for exam 190924 problem 5:"
$ws.Range("K17").Value = "E(Z|y=2)"
$ws.Range("L17").Formula = "=X6"
$ws.Range("K19").Value = "And you can add more rows and columns in general if a problem wants it:"
$ws.Range("K20").Value = "Still works!"

# Update the selected cell/range to match the final saved selection
$ws.Range("K23").Select()
